$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values 5 in C21:F21 (style already set to s="2", only values change)
$ws.Range("C21:F21").Value = 5

# Update the active cell / selection to G21 as recorded in the sheet view
$ws.Range("G21").Select()
